$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 28
$ws.Range("H28").Value = 1057.5454
$ws.Range("I28").Value = 847.3333
$ws.Range("K28").Value = 847.3333
$ws.Range("M28").Value = -362.3333
# Row 40
$ws.Range("H40").Value = 2662.3713
$ws.Range("J40").Value = 3465.8572
$ws.Range("L40").Value = 3465.8572
$ws.Range("N40").Value = -3815.8572
# Row 74
$ws.Range("H74").Value = 41753880
$ws.Range("I74").Value = 47717150
$ws.Range("K74").Value = 47717150
$ws.Range("M74").Value = -47716214
# Row 77
$ws.Range("H77").Value = 41753880
$ws.Range("I77").Value = 47717150
$ws.Range("K77").Value = 238585750
$ws.Range("M77").Value = -238581070
# Row 80
$ws.Range("H80").Value = 866.125
$ws.Range("J80").Value = 1062.5
$ws.Range("L80").Value = 3187.5
$ws.Range("N80").Value = -5183.5
# Row 83
$ws.Range("H83").Value = 866.125
$ws.Range("J83").Value = 1062.5
$ws.Range("L83").Value = 9562.5
$ws.Range("N83").Value = -19546.5
# Row 107
$ws.Range("H107").Value = 1300
$ws.Range("I107").Value = 1601
$ws.Range("J107").Value = 923.75
$ws.Range("K107").Value = 1601
$ws.Range("L107").Value = 923.75
$ws.Range("M107").Value = 319
$ws.Range("N107").Value = -4763.75
# Row 112
$ws.Range("H112").Value = 79843
$ws.Range("I112").Value = 85348.25
$ws.Range("J112").Value = 75124.21000000001
$ws.Range("K112").Value = 256044.75
$ws.Range("L112").Value = 225372.63
$ws.Range("M112").Value = -254936.75
$ws.Range("N112").Value = -227588.63
# Row 127
$ws.Range("H127").Value = 2227.4443
$ws.Range("I127").Value = 2206.7144
$ws.Range("K127").Value = 6620.1432
$ws.Range("M127").Value = -1660.1432
# Row 129
$ws.Range("H129").Value = 1801.1177
$ws.Range("J129").Value = 1725.091
$ws.Range("L129").Value = 5175.272999999999
$ws.Range("N129").Value = -15175.273
# Row 131
$ws.Range("H131").Value = 1019
$ws.Range("I131").Value = 1019
$ws.Range("K131").Value = 3057
$ws.Range("M131").Value = 1983
# Row 132
$ws.Range("H132").Value = 1681.12
$ws.Range("I132").Value = 1717.25
$ws.Range("K132").Value = 5151.75
$ws.Range("M132").Value = -2621.75
# Row 137
$ws.Range("H137").Value = 1974.7587
$ws.Range("I137").Value = 2022.909
$ws.Range("K137").Value = 6068.727000000001
$ws.Range("M137").Value = -3518.727000000001
# Row 141
$ws.Range("H141").Value = 3114
$ws.Range("J141").Value = 4786.5
$ws.Range("L141").Value = 14359.5
$ws.Range("N141").Value = -24719.5

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 5
$ws.Range("H5").Value = 431.5
$ws.Range("I5").Value = 174.5
$ws.Range("K5").Value = 174.5
$ws.Range("M5").Value = -62.5
# Row 45
$ws.Range("H45").Value = 1361
$ws.Range("I45").Value = 1243.75
$ws.Range("K45").Value = 1243.75
$ws.Range("M45").Value = -866.75
# Row 63
$ws.Range("H63").Value = 3797.2
$ws.Range("I63").Value = 3797.2
$ws.Range("K63").Value = 3797.2
$ws.Range("M63").Value = -3111.2
# Row 66
$ws.Range("H66").Value = 3797.2
$ws.Range("I66").Value = 3797.2
$ws.Range("K66").Value = 18986
$ws.Range("M66").Value = -15554
# Row 88
$ws.Range("H88").Value = 3248.75
$ws.Range("I88").Value = 3000
$ws.Range("J88").Value = 3331.6667
$ws.Range("K88").Value = 3000
$ws.Range("L88").Value = 3331.6667
$ws.Range("M88").Value = -2594
$ws.Range("N88").Value = -4143.6667
# Row 91
$ws.Range("H91").Value = 3248.75
$ws.Range("I91").Value = 3000
$ws.Range("J91").Value = 3331.6667
$ws.Range("K91").Value = 3000
$ws.Range("L91").Value = 3331.6667
$ws.Range("M91").Value = -1596
$ws.Range("N91").Value = -6139.6667
# Row 110
$ws.Range("H110").Value = 37435.137
$ws.Range("I110").Value = 44401.332
$ws.Range("K110").Value = 44401.332
$ws.Range("M110").Value = -42356.332

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 4
$ws.Range("H4").Value = 431.5
$ws.Range("I4").Value = 174.5
$ws.Range("K4").Value = 174.5
$ws.Range("M4").Value = -59.5
# Row 86
$ws.Range("H86").Value = 4881.8335
$ws.Range("I86").Value = 4998.5
$ws.Range("J86").Value = 4648.5
$ws.Range("K86").Value = 4998.5
$ws.Range("L86").Value = 4648.5
$ws.Range("M86").Value = -3875.5
$ws.Range("N86").Value = -6894.5
# Row 89
$ws.Range("H89").Value = 4881.8335
$ws.Range("I89").Value = 4998.5
$ws.Range("J89").Value = 4648.5
$ws.Range("K89").Value = 24992.5
$ws.Range("L89").Value = 23242.5
$ws.Range("M89").Value = -19376.5
$ws.Range("N89").Value = -34474.5
# Row 94
$ws.Range("H94").Value = 961.86957
$ws.Range("I94").Value = 948.5333000000001
$ws.Range("J94").Value = 986.875
$ws.Range("K94").Value = 948.5333000000001
$ws.Range("L94").Value = 986.875
$ws.Range("M94").Value = -497.5333000000001
$ws.Range("N94").Value = -1888.875
# Row 134
$ws.Range("H134").Value = 31251652
$ws.Range("I134").Value = 31251652
$ws.Range("K134").Value = 93754956
$ws.Range("M134").Value = -93752421

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 116
$ws.Range("H116").Value = 67500
$ws.Range("J116").Value = 67500
$ws.Range("L116").Value = 67500
$ws.Range("N116").Value = -76678
# Row 134
$ws.Range("H134").Value = 27779962
$ws.Range("I134").Value = 35715950
$ws.Range("K134").Value = 107147850
$ws.Range("M134").Value = -107145315

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 129
$ws.Range("H129").Value = 1663
$ws.Range("I129").Value = 736.75
$ws.Range("J129").Value = 4133
$ws.Range("K129").Value = 2210.25
$ws.Range("L129").Value = 12399
$ws.Range("M129").Value = 2789.75
$ws.Range("N129").Value = -22399
# Row 131
$ws.Range("H131").Value = 1653.75
$ws.Range("J131").Value = 1916
$ws.Range("L131").Value = 5748
$ws.Range("N131").Value = -15828
# Row 139
$ws.Range("H139").Value = 2986
$ws.Range("I139").Value = 2986
$ws.Range("K139").Value = 8958
$ws.Range("M139").Value = -3818
# Row 140
$ws.Range("H140").Value = 1482.3549
$ws.Range("I140").Value = 755.5357
$ws.Range("K140").Value = 2266.6071
$ws.Range("M140").Value = 2913.3929

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 6640.6665
$ws.Range("I70").Value = 9939.799999999999
$ws.Range("J70").Value = 4284.143
$ws.Range("K70").Value = 9939.799999999999
$ws.Range("L70").Value = 4284.143
$ws.Range("M70").Value = -9669.799999999999
$ws.Range("N70").Value = -4824.143
# Row 73
$ws.Range("H73").Value = 6640.6665
$ws.Range("I73").Value = 9939.799999999999
$ws.Range("J73").Value = 4284.143
$ws.Range("K73").Value = 9939.799999999999
$ws.Range("L73").Value = 4284.143
$ws.Range("M73").Value = -9003.799999999999
$ws.Range("N73").Value = -6156.143
# Row 80
$ws.Range("H80").Value = 2624.75
$ws.Range("I80").Value = 2624.75
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 2624.75
$ws.Range("L80").Value = 0
$ws.Range("M80").Value = -1626.75
$ws.Range("N80").Value = $null
# Row 83
$ws.Range("H83").Value = 2624.75
$ws.Range("I83").Value = 2624.75
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 13123.75
$ws.Range("L83").Value = 0
$ws.Range("M83").Value = -8131.75
$ws.Range("N83").Value = $null
# Row 113
$ws.Range("H113").Value = 79392.46000000001
$ws.Range("I113").Value = 85258.5
$ws.Range("J113").Value = 9000
$ws.Range("K113").Value = 85258.5
$ws.Range("L113").Value = 9000
$ws.Range("M113").Value = -83088.5
$ws.Range("N113").Value = -13340
# Row 132
$ws.Range("H132").Value = 11370830
$ws.Range("I132").Value = 13892000
$ws.Range("J132").Value = 25562.5
$ws.Range("K132").Value = 41676000
$ws.Range("L132").Value = 76687.5
$ws.Range("M132").Value = -41673470
$ws.Range("N132").Value = -81747.5

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Range("H16").Value = 1232.36
$ws.Range("I16").Value = 292.58823
$ws.Range("J16").Value = 3229.375
$ws.Range("K16").Value = 292.58823
$ws.Range("L16").Value = 3229.375
$ws.Range("M16").Value = -122.58823
$ws.Range("N16").Value = -3569.375
# Row 32
$ws.Range("H32").Value = 2115.75
$ws.Range("I32").Value = 2115.75
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 2115.75
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -1798.75
$ws.Range("N32").Value = $null
# Row 61
$ws.Range("H61").Value = 3153.2307
$ws.Range("I61").Value = 2916
$ws.Range("K61").Value = 2916
$ws.Range("M61").Value = -2714
# Row 94
$ws.Range("H94").Value = 25000
$ws.Range("J94").Value = 25000
$ws.Range("L94").Value = 25000
$ws.Range("N94").Value = -26352
# Row 113
$ws.Range("H113").Value = 3153.2307
$ws.Range("I113").Value = 2916
$ws.Range("K113").Value = 2916
$ws.Range("M113").Value = -746
# Row 122
$ws.Range("H122").Value = 4544.893
$ws.Range("I122").Value = 4321.423
$ws.Range("K122").Value = 12964.269
$ws.Range("M122").Value = -10514.269

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 41
$ws.Range("H41").Value = 49999.5
$ws.Range("J41").Value = 49999.5
$ws.Range("L41").Value = 49999.5
$ws.Range("N41").Value = -50779.5
# Row 62
$ws.Range("H62").Value = 10000
$ws.Range("I62").Value = 10000
$ws.Range("K62").Value = 10000
$ws.Range("M62").Value = -9376
# Row 65
$ws.Range("H65").Value = 10000
$ws.Range("I65").Value = 10000
$ws.Range("K65").Value = 50000
$ws.Range("M65").Value = -46880
# Row 81
$ws.Range("H81").Value = 3286.5557
$ws.Range("I81").Value = 3134.875
$ws.Range("J81").Value = 4500
$ws.Range("K81").Value = 6269.75
$ws.Range("L81").Value = 9000
$ws.Range("M81").Value = -5208.75
$ws.Range("N81").Value = -11122
# Row 84
$ws.Range("H84").Value = 3286.5557
$ws.Range("I84").Value = 3134.875
$ws.Range("J84").Value = 4500
$ws.Range("K84").Value = 31348.75
$ws.Range("L84").Value = 45000
$ws.Range("M84").Value = -26044.75
$ws.Range("N84").Value = -55608

